# Update cryptos list cell values (prices + 1h volume %) to match latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.092.35'
$ws.Range('E2').Value = '  +2.32%  '

$ws.Range('D3').Value = '3.469.18'
$ws.Range('E3').Value = '  +2.16%  '

$ws.Range('E4').Value = '  -0.08%  '

$ws.Range('D5').Value = '''578.10'
$ws.Range('E5').Value = '  +0.27%  '

$ws.Range('D6').Value = '''148.22'
$ws.Range('E6').Value = '  +3.31%  '

$ws.Range('D7').Value = '3.465.07'
$ws.Range('E7').Value = '  +1.99%  '

$ws.Range('E8').Value = '  +0.01%  '

$ws.Range('D9').Value = '''0.479'
$ws.Range('E9').Value = '  +1.29%  '

$ws.Range('E10').Value = '  +0.49%  '

$ws.Range('E11').Value = '  +1.78%  '

$ws.Range('D12').Value = '''0.404'
$ws.Range('E12').Value = '  +4.64%  '

$ws.Range('D13').Value = '4.060.83'
$ws.Range('E13').Value = '  +2.05%  '

$ws.Range('D14').Value = '''29.90'
$ws.Range('E14').Value = '  +6.51%  '

$ws.Range('E15').Value = '  +2.65%  '

$ws.Range('D16').Value = '3.479.06'
$ws.Range('E16').Value = '  +2.40%  '

$ws.Range('E17').Value = '  +0.54%  '

$ws.Range('D18').Value = '63.022.70'
$ws.Range('E18').Value = '  +2.12%  '

$ws.Range('D19').Value = '''6.35'
$ws.Range('E19').Value = '  +3.47%  '

$ws.Range('D20').Value = '''14.40'
$ws.Range('E20').Value = '  +5.40%  '

$ws.Range('D21').Value = '''9.27'
$ws.Range('E21').Value = '  +1.26%  '

$ws.Range('D22').Value = '''388.76'
$ws.Range('E22').Value = '  +0.15%  '

$ws.Range('E23').Value = '  +1.60%  '

$ws.Range('D24').Value = '''74.70'
$ws.Range('E24').Value = '  +0.16%  '

$ws.Range('E25').Value = '  +0.02%  '

$ws.Range('D26').Value = '3.610.14'

$ws.Range('D27').Value = '''0.0000116'
$ws.Range('E27').Value = '  +0.85%  '

$ws.Range('E28').Value = '  -1.97%  '

$ws.Range('D29').Value = '''7.63'
$ws.Range('E29').Value = '  +3.20%  '

$ws.Range('D30').Value = '''0.999'
$ws.Range('E30').Value = '  +0.19%  '

$ws.Range('D31').Value = '''8.17'
$ws.Range('E31').Value = '  +2.17%  '

$ws.Range('E32').Value = '  -0.80%  '

$ws.Range('E33').Value = '  +0.09%  '

$ws.Range('B34').Value = 'Fetch.AI'
$ws.Range('C34').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D34').Value = '''1.38'
$ws.Range('E34').Value = '  -1.87%  '

$ws.Range('B35').Value = 'EthereumClassic'
$ws.Range('C35').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D35').Value = '''23.68'
$ws.Range('E35').Value = '  +1.32%  '

$ws.Range('B36').Value = 'NEARProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D36').Value = '''5.31'
$ws.Range('E36').Value = '  +4.22%  '

$ws.Range('B37').Value = 'Aptos'
$ws.Range('C37').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D37').Value = '''7.08'
$ws.Range('E37').Value = '  +2.08%  '

$ws.Range('B38').Value = 'EnergySwap'
$ws.Range('C38').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D38').Value = '''31.90'
$ws.Range('E38').Value = '  +18.11%  '

$ws.Range('D39').Value = '''169.75'
$ws.Range('E39').Value = '  +0.42%  '

$ws.Range('D40').Value = '''1.57'
$ws.Range('E40').Value = '  +6.03%  '

$ws.Range('D41').Value = '3.505.80'
$ws.Range('E41').Value = '  +2.19%  '

$ws.Range('D42').Value = '''0.0758'
$ws.Range('E42').Value = '  -0.42%  '

$ws.Range('D43').Value = '''0.798'
$ws.Range('E43').Value = '  +1.78%  '

$ws.Range('D44').Value = '''42.38'
$ws.Range('E44').Value = '  -0.10%  '

$ws.Range('B45').Value = 'ONDO'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D45').Value = '''1.22'
$ws.Range('E45').Value = '  +5.09%  '

$ws.Range('B46').Value = 'Filecoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D46').Value = '''4.47'
$ws.Range('E46').Value = '  +0.90%  '

$ws.Range('B47').Value = 'Stacks'
$ws.Range('C47').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D47').Value = '''1.72'
$ws.Range('E47').Value = '  +3.21%  '

$ws.Range('D48').Value = '2.619.69'
$ws.Range('E48').Value = '  +5.65%  '

$ws.Range('D49').Value = '''2.29'
$ws.Range('E49').Value = '  +12.04%  '

$ws.Range('D50').Value = '''23.02'
$ws.Range('E50').Value = '  +0.99%  '

$ws.Range('D51').Value = '''6.75'
$ws.Range('E51').Value = '  +1.59%  '
